$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Extend the sentence "...is specified as well to include a
#    percentage sign." with the new explanatory sentence about
#    miss_decimal(#) and su_decimal(#). Built up run-by-run so the
#    VerbatimChar character style lands only on the option names.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("percentage sign.") | Out-Null
# Extend the existing run's text in place (stays one <w:r>), rather than
# inserting a new run after it.
$r.Text = "percentage sign. When denominators or missing data summaries are included in the table the options"
$r.Collapse(0)

$r.InsertAfter(" ")
$r.Collapse(0)

$r.InsertAfter("miss_decimal(#)")
$r.Style = "VerbatimChar"
$r.Collapse(0)

$r.InsertAfter(" ")
$r.Collapse(0)

$r.InsertAfter("and")
$r.Collapse(0)

$r.InsertAfter(" ")
$r.Collapse(0)

$r.InsertAfter("su_decimal(#)")
$r.Style = "VerbatimChar"
$r.Collapse(0)

$r.InsertAfter(" ")
$r.Collapse(0)

$r.InsertAfter("can be used to independently control the number of decimal places reported for summary statistics and the percent of missing/nonmissing observations.")
$r.Collapse(0)

# ---------------------------------------------------------------------
# 2. Update the "post" example line: the two-group column headers
#    change from N-based labels to Missing-based labels. Assigning
#    .Text directly (rather than Find.Execute's Replacement.Text)
#    avoids the straight quotes being "smart-quoted" by autoformat.
# ---------------------------------------------------------------------
$rPost1 = $d.Content
$rPost1.Find.Execute('("N 1")') | Out-Null
$rPost1.Text = '("Missing 1")'

$rPost2 = $d.Content
$rPost2.Find.Execute('("N 2")') | Out-Null
$rPost2.Text = '("Missing 0")'

$rPost3 = $d.Content
$rPost3.Find.Execute('("Summary 2")') | Out-Null
$rPost3.Text = '("Summary 0")'

# ---------------------------------------------------------------------
# 3. Append the new miss_decimal()/su_decimal()/decimal() options to
#    the "pt_base age" and "pt_base qol" example lines.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
  "over_grps(1, 0) type(cont) su_label(append) cat_col  missing(cols cond %) order(group_over) per", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "over_grps(1, 0) type(cont) su_label(append) cat_col  missing(cols cond %) order(group_over) per  miss_decimal(2) su_decimal(0)", `
  2) | Out-Null

$d.Content.Find.Execute( `
  "over_grps(1, 0) type(skew) su_label(append) cat_col  missing(cols  cond %) order(group_over) per", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  "over_grps(1, 0) type(skew) su_label(append) cat_col  missing(cols  cond %) order(group_over) per  miss_decimal(2) decimal(1)", `
  2) | Out-Null
